$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "50.205.46"
$ws.Range("E2").Value = "  +4.42%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.644.02"
$ws.Range("E3").Value = "  +5.85%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.68"
$ws.Range("E5").Value = "  +2.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "111.59"
$ws.Range("E6").Value = "  +3.67%  "
$ws.Range("E7").Value = "  +1.75%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.563"
$ws.Range("E9").Value = "  +4.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.95"
$ws.Range("E10").Value = "  +3.40%  "
$ws.Range("E11").Value = "  +2.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0822"
$ws.Range("E12").Value = "  +1.06%  "
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.31"
$ws.Range("E14").Value = "  +2.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.048.51"
$ws.Range("E15").Value = "  +5.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.646.96"
$ws.Range("E16").Value = "  +6.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.880"
$ws.Range("E17").Value = "  +5.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "50.047.19"
$ws.Range("E18").Value = "  +4.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.08"
$ws.Range("E19").Value = "  +12.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.41"
$ws.Range("E20").Value = "  +3.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.84"
$ws.Range("E21").Value = "  +1.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0963"
$ws.Range("E22").Value = "  +2.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.10"
$ws.Range("E23").Value = "  +2.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "279.70"
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.61"
$ws.Range("E25").Value = "  +3.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.70"
$ws.Range("E26").Value = "  +4.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.93"
$ws.Range("E28").Value = "  +5.39%  "
$ws.Range("E29").Value = "  +6.40%  "
$ws.Range("E30").Value = "  +2.17%  "
$ws.Range("E31").Value = "  +3.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.85"
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.87"
$ws.Range("E33").Value = "  +1.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.46"
$ws.Range("E34").Value = "  +2.94%  "
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0795"
$ws.Range("E36").Value = "  +1.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.07"
$ws.Range("E37").Value = "  +6.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.77"
$ws.Range("E38").Value = "  +3.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.11"
$ws.Range("E39").Value = "  +7.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.113"
$ws.Range("E40").Value = "  +1.44%  "
$ws.Range("E41").Value = "  +1.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.64"
$ws.Range("E42").Value = "  +6.10%  "
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0316"
$ws.Range("E44").Value = "  +5.10%  "
$ws.Range("E45").Value = "  +7.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.063.49"
$ws.Range("E46").Value = "  +2.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.32"
$ws.Range("E47").Value = "  +16.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.01"
$ws.Range("E48").Value = "  +8.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.05"
$ws.Range("E49").Value = "  +0.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.39"
$ws.Range("E50").Value = "  +4.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.87"
$ws.Range("E51").Value = "  +2.06%  "
